$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume figures.
# NumberFormat="@" + Style="Normal" keeps numeric-looking strings (e.g. "548.08")
# stored as TEXT (matches original inlineStr cells) instead of being coerced to a number,
# while resetting the cell style back to the default (no stray "s" attribute).

# Row 2
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '62.827.87'
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '  +1.96%  '
$r.Style = "Normal"

# Row 3
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '3.038.13'
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '  +0.37%  '
$r.Style = "Normal"

# Row 4
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '  -0.13%  '
$r.Style = "Normal"

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '548.08'
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '  +2.29%  '
$r.Style = "Normal"

# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '137.32'
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '  +0.93%  '
$r.Style = "Normal"

# Row 7
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '  -0.12%  '
$r.Style = "Normal"

# Row 8
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '3.033.84'
$r.Style = "Normal"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '  +0.33%  '
$r.Style = "Normal"

# Row 9
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '  -0.23%  '
$r.Style = "Normal"

# Row 10
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '6.29'
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '  +2.21%  '
$r.Style = "Normal"

# Row 11
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.149'
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '  -1.04%  '
$r.Style = "Normal"

# Row 12
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '  -0.78%  '
$r.Style = "Normal"

# Row 13
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '  +1.31%  '
$r.Style = "Normal"

# Row 14
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '34.33'
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '  -1.19%  '
$r.Style = "Normal"

# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '3.534.60'
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '  +0.23%  '
$r.Style = "Normal"

# Row 16
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '62.938.82'
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '  +1.90%  '
$r.Style = "Normal"

# Row 17
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '3.034.71'
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '  +0.01%  '
$r.Style = "Normal"

# Row 18
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '  -2.30%  '
$r.Style = "Normal"

# Row 19
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '6.67'
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = '  -0.32%  '
$r.Style = "Normal"

# Row 20
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '478.03'
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '  +1.54%  '
$r.Style = "Normal"

# Row 21
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '13.55'
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '  +1.12%  '
$r.Style = "Normal"

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '0.667'
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '  -2.60%  '
$r.Style = "Normal"

# Row 23
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '7.14'
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '  +1.45%  '
$r.Style = "Normal"

# Row 24
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '80.14'
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '  +0.26%  '
$r.Style = "Normal"

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '12.33'
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '  +0.76%  '
$r.Style = "Normal"

# Row 26
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '  -0.02%  '
$r.Style = "Normal"

# Row 27
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = '  +1.13%  '
$r.Style = "Normal"

# Row 28
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '7.81'
$r.Style = "Normal"
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = '  -1.40%  '
$r.Style = "Normal"

# Row 29
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = '  -0.37%  '
$r.Style = "Normal"

# Row 30
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = '  +1.85%  '
$r.Style = "Normal"

# Row 31
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '25.76'
$r.Style = "Normal"
$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = '  -0.37%  '
$r.Style = "Normal"

# Row 32
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = '  -0.35%  '
$r.Style = "Normal"

# Row 33
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = '  +3.64%  '
$r.Style = "Normal"

# Row 34
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '5.64'
$r.Style = "Normal"
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = '  +1.48%  '
$r.Style = "Normal"

# Row 35
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '55.18'
$r.Style = "Normal"
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = '  -0.89%  '
$r.Style = "Normal"

# Row 36
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = '  -1.08%  '
$r.Style = "Normal"

# Row 37
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '459.28'
$r.Style = "Normal"
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = '  -2.73%  '
$r.Style = "Normal"

# Row 38
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.0807'
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = '  +1.17%  '
$r.Style = "Normal"

# Row 39
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '  -0.06%  '
$r.Style = "Normal"

# Row 40
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '3.056.80'
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '  -5.29%  '
$r.Style = "Normal"

# Row 41
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '  -1.05%  '
$r.Style = "Normal"

# Row 42
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '8.19'
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '  -0.36%  '
$r.Style = "Normal"

# Row 43
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '2.53'
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '  -0.18%  '
$r.Style = "Normal"

# Row 44
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '27.88'
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '  +0.11%  '
$r.Style = "Normal"

# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.250'
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '  -0.45%  '
$r.Style = "Normal"

# Row 46
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '  -0.11%  '
$r.Style = "Normal"

# Row 47
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '  -0.61%  '
$r.Style = "Normal"

# Row 48
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '  -0.56%  '
$r.Style = "Normal"

# Row 49
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '116.10'
$r.Style = "Normal"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '  -3.08%  '
$r.Style = "Normal"

# Row 50
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '  -0.24%  '
$r.Style = "Normal"

# Row 51
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '2.04'
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '  +0.53%  '
$r.Style = "Normal"
